# Updated cryptos list on Fri Mar 31 17:56:37 UTC 2023 with GitHub Actions
#
# Refresh the Price (D) and Volume(1h) (E) columns for every coin row on
# Sheet1, plus the rank-17/18 swap between Litecoin and ShibaInu (name,
# link, price and volume all move together).
#
# Several Price values look like plain numbers (e.g. "317.53"); setting
# .Value directly on those would let Excel coerce them to a numeric type
# (and round-trip through floating point), so those cells get their
# NumberFormat forced to text ("@") first, exactly like Excel does when a
# user types a numeric-looking value into a pre-formatted text cell.
# Cells whose new text is unambiguous (contains two dots, % sign, letters,
# spaces, etc.) are left alone since Excel already keeps those as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.331.28'
$ws.Range("E2").Value = '  +1.55%  '

$ws.Range("D3").Value = '1.825.62'
$ws.Range("E3").Value = '  +2.82%  '

$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.53'
$ws.Range("E5").Value = '  +0.77%  '

$ws.Range("E6").Value = '  -0.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5331'
$ws.Range("E7").Value = '  -0.28%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4081'
$ws.Range("E8").Value = '  +9.73%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07599'
$ws.Range("E9").Value = '  +2.92%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.86'
$ws.Range("E10").Value = '  +1.14%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.102'
$ws.Range("E11").Value = '  +1.24%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.324'
$ws.Range("E12").Value = '  +4.42%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.001'
$ws.Range("E13").Value = '  -0.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.559'
$ws.Range("E14").Value = '  +5.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '20.75'
$ws.Range("E15").Value = '  +1.58%  '

$ws.Range("D16").Value = '1.828.76'
$ws.Range("E16").Value = '  +2.75%  '

# Rank 17 / 18 swap: Litecoin <-> ShibaInu (name, link, price, volume)
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001072'
$ws.Range("E17").Value = '  +2.29%  '

$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '89.21'
$ws.Range("E18").Value = '  +2.04%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06616'
$ws.Range("E19").Value = '  +2.75%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.58'
$ws.Range("E20").Value = '  +1.56%  '

$ws.Range("E21").Value = '  -0.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.087'
$ws.Range("E22").Value = '  +3.62%  '

$ws.Range("D23").Value = '28.379.21'
$ws.Range("E23").Value = '  +1.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.16'
$ws.Range("E24").Value = '  +0.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.181'
$ws.Range("E25").Value = '  +4.16%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.458'
$ws.Range("E26").Value = '  +8.37%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '156.77'
$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.53'
$ws.Range("E28").Value = '  +2.16%  '

$ws.Range("D29").Value = '2.039.74'
$ws.Range("E29").Value = '  +3.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.55'
$ws.Range("E30").Value = '  +3.25%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.116'
$ws.Range("E31").Value = '  +1.45%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1092'
$ws.Range("E32").Value = '  +5.48%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.657'
$ws.Range("E33").Value = '  +3.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.629'
$ws.Range("E34").Value = '  -0.53%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07180'
$ws.Range("E35").Value = '  +13.17%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.2254'
$ws.Range("E36").Value = '  +1.30%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02336'
$ws.Range("E37").Value = '  +3.43%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.205'
$ws.Range("E38").Value = '  +5.14%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.839'
$ws.Range("E39").Value = '  +5.51%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6266'
$ws.Range("E40").Value = '  +2.27%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.27'
$ws.Range("E41").Value = '  +3.01%  '

$ws.Range("E42").Value = '  +1.35%  '

$ws.Range("E43").Value = '  -0.16%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.403'
$ws.Range("E44").Value = '  -2.22%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.48'
$ws.Range("E45").Value = '  +2.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.707'
$ws.Range("E46").Value = '  +1.32%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5835'
$ws.Range("E47").Value = '  +1.91%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.58'
$ws.Range("E48").Value = '  +0.45%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.984'
$ws.Range("E49").Value = '  +3.32%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.210'
$ws.Range("E50").Value = '  +1.68%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06888'
$ws.Range("E51").Value = '  +1.12%  '
